$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, shifting existing rows 130:240 down to 131:241
$ws.Rows("130:130").Insert()

# Populate the newly inserted row 130 with the new weekly record
$ws.Cells.Item(130, 1).Value = 3
$ws.Cells.Item(130, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(130, 3).Value = "Coquimbo"
$ws.Cells.Item(130, 4).Value = 44957
$ws.Cells.Item(130, 5).Value = 5
$ws.Cells.Item(130, 6).Value = 100112030
$ws.Cells.Item(130, 7).Value = "Poroto granado"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 73
$ws.Cells.Item(130, 11).Value = 41000
$ws.Cells.Item(130, 12).Value = 42000
$ws.Cells.Item(130, 13).Value = 41521
$ws.Cells.Item(130, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(130, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(130, 16).Value = 1661
$ws.Cells.Item(130, 17).Value = 25
$ws.Cells.Item(130, 18).Value = "Hortaliza"
